# Apply the LOM3111.xlsx worksheet restructuring described in the commit diff.
# The sheet rows 10-23 are reshuffled (headers/values shift positions) and
# trailing rows 24-25 are removed, shrinking the sheet from A1:C25 to A1:C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that must become empty in the target layout ---
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()

# --- Overwrite cells with their target content ---
$ws.Range("B10").Value = '5983729 - Fernando Vernilli Junior'
$ws.Range("C10").Value = '5983729 - Fernando Vernilli Junior'
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = '01/01/2022'
$ws.Range("C13").Value = '01/01/2022'
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '5983729 - Fernando Vernilli Junior'
$ws.Range("C15").Value = '5983729 - Fernando Vernilli Junior'
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("C18").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Avaliação individual do comportamento do aluno frente aos trabalhos experimentais (AC), Relatórios sobre os testes experimentais (MAR) e prova experimental final (PE).'
$ws.Range("C19").Value = 'Avaliação individual do comportamento do aluno frente aos trabalhos experimentais (AC), Relatórios sobre os testes experimentais (MAR) e prova experimental final (PE).'
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = '1 - Média aritmética das notas dos relatórios, com peso 1 (MAR), 2 – avaliação comportamental, peso 1 (AC) e uma prova experimental (PE), no final do semestre letivo, com peso 2.A nota final (NF) será calculada pela equação (MAR+AC+2PE)/4. NF igual ou superior a 5: aprovação.'
$ws.Range("C20").Value = '1 - Média aritmética das notas dos relatórios, com peso 1 (MAR), 2 – avaliação comportamental, peso 1 (AC) e uma prova experimental (PE), no final do semestre letivo, com peso 2.A nota final (NF) será calculada pela equação (MAR+AC+2PE)/4. NF igual ou superior a 5: aprovação.'
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Devido à característica da disciplina não'
$ws.Range("C21").Value = 'Devido à característica da disciplina não'
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B23").Value = 'LOM3073 -  Processamento de Cerâmicas  (Requisito fraco)' + [char]10 + ''
$ws.Range("C23").Value = 'LOM3073 -  Processamento de Cerâmicas  (Requisito fraco)' + [char]10 + ''

# --- Adjust row heights to match the target layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30

# --- Remove the now-obsolete trailing rows (old Requisitos value row + LOM3073 row) ---
$ws.Rows.Item(24).Resize(2).Delete()

Write-Host ("Final UsedRange: " + $ws.UsedRange.Address())
